# Apply the cryptos price/volume update described by the commit:
# "Updated cryptos list on Sat Jul 22 03:27:50 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for Price cells whose new value would otherwise
# be auto-recognised by Excel as a number (e.g. "244.76", "1.000"), which
# would silently drop meaningful trailing zeros / change the stored type.
$textPriceRows = @(5,6,7,8,9,10,11,12,13,15,16,18,19,20,22,23,25,26,27,28,29,30,32,33,34,35,36,37,38,39,40,41,42,43,45,47,48,49,50,51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = '29.939.78'
$ws.Range("E2").Value = '  +0.07%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.895.87'
$ws.Range("E3").Value = '  -0.17%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.02%  '

# Row 5 - XRP
$ws.Range("D5").Value = '0.7775'
$ws.Range("E5").Value = '  -2.29%  '

# Row 6 - BNB
$ws.Range("D6").Value = '244.76'
$ws.Range("E6").Value = '  +0.10%  '

# Row 7 - USDC
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.01%  '

# Row 8 - Cardano
$ws.Range("D8").Value = '0.3142'
$ws.Range("E8").Value = '  -1.06%  '

# Row 9 - Solana
$ws.Range("D9").Value = '25.96'
$ws.Range("E9").Value = '  +1.61%  '

# Row 10 - Dogecoin
$ws.Range("D10").Value = '0.07267'
$ws.Range("E10").Value = '  +1.08%  '

# Row 11 - TRON
$ws.Range("D11").Value = '0.09313'
$ws.Range("E11").Value = '  +14.73%  '

# Row 12 - Polygon
$ws.Range("D12").Value = '0.7748'
$ws.Range("E12").Value = '  +0.37%  '

# Row 13 - Polkadot (was WrappedEther)
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.488'
$ws.Range("E13").Value = '  -2.92%  '

# Row 14 - WrappedEther (was Polkadot)
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.912.08'
$ws.Range("E14").Value = '  +1.41%  '

# Row 15 - Litecoin
$ws.Range("D15").Value = '94.80'
$ws.Range("E15").Value = '  +2.10%  '

# Row 16 - Uniswap
$ws.Range("D16").Value = '6.253'
$ws.Range("E16").Value = '  +1.00%  '

# Row 17 - WrappedBTC
$ws.Range("D17").Value = '30.015.49'
$ws.Range("E17").Value = '  +0.37%  '

# Row 18 - Avalanche
$ws.Range("D18").Value = '14.01'
$ws.Range("E18").Value = '  +0.06%  '

# Row 19 - BitcoinCash
$ws.Range("D19").Value = '247.22'
$ws.Range("E19").Value = '  +0.57%  '

# Row 20 - ShibaInu
$ws.Range("D20").Value = '0.000007887'
$ws.Range("E20").Value = '  +1.33%  '

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = '2.178.75'
$ws.Range("E21").Value = '  +2.65%  '

# Row 22 - Chainlink
$ws.Range("D22").Value = '8.195'
$ws.Range("E22").Value = '  -0.45%  '

# Row 23 - Dai
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.01%  '

# Row 24 - BinanceUSD
$ws.Range("E24").Value = '  -0.04%  '

# Row 25 - Stellar
$ws.Range("D25").Value = '0.1597'
$ws.Range("E25").Value = '  -4.50%  '

# Row 26 - Cosmos
$ws.Range("D26").Value = '9.554'
$ws.Range("E26").Value = '  +0.56%  '

# Row 27 - Monero
$ws.Range("D27").Value = '162.32'
$ws.Range("E27").Value = '  -1.22%  '

# Row 28 - EthereumClassic
$ws.Range("D28").Value = '18.83'
$ws.Range("E28").Value = '  +0.26%  '

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = '2.053'
$ws.Range("E29").Value = '  -1.48%  '

# Row 30 - Toncoin
$ws.Range("D30").Value = '1.425'
$ws.Range("E30").Value = '  +1.04%  '

# Row 31 - PancakeSwap
$ws.Range("E31").Value = '  +0.18%  '

# Row 32 - Filecoin
$ws.Range("D32").Value = '4.546'
$ws.Range("E32").Value = '  +1.09%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = '4.125'
$ws.Range("E33").Value = '  +0.67%  '

# Row 34 - Hedera
$ws.Range("D34").Value = '0.05514'
$ws.Range("E34").Value = '  -2.21%  '

# Row 35 - ARBITRUM
$ws.Range("D35").Value = '1.251'
$ws.Range("E35").Value = '  -3.25%  '

# Row 36 - ImmutableX
$ws.Range("D36").Value = '0.7583'
$ws.Range("E36").Value = '  +1.44%  '

# Row 37 - Frax
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  +0.42%  '

# Row 38 - HuobiToken
$ws.Range("D38").Value = '2.692'
$ws.Range("E38").Value = '  +2.13%  '

# Row 39 - VeChain
$ws.Range("D39").Value = '0.01974'
$ws.Range("E39").Value = '  +1.80%  '

# Row 40 - MXToken
$ws.Range("D40").Value = '2.792'
$ws.Range("E40").Value = '  +0.36%  '

# Row 41 - TheSandbox
$ws.Range("D41").Value = '0.4525'
$ws.Range("E41").Value = '  +1.94%  '

# Row 42 - Aave
$ws.Range("D42").Value = '74.42'
$ws.Range("E42").Value = '  -0.61%  '

# Row 43 - FraxShare
$ws.Range("D43").Value = '6.094'
$ws.Range("E43").Value = '  +2.13%  '

# Row 44 - Maker
$ws.Range("D44").Value = '1.093.83'
$ws.Range("E44").Value = '  -6.22%  '

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = '0.8534'

# Row 46 - PaxDollar
$ws.Range("E46").Value = '  -0.01%  '

# Row 47 - RenderToken
$ws.Range("D47").Value = '1.899'
$ws.Range("E47").Value = '  +0.44%  '

# Row 48 - Quant
$ws.Range("D48").Value = '102.89'

# Row 49 - Aptos
$ws.Range("D49").Value = '7.620'
$ws.Range("E49").Value = '  +1.50%  '

# Row 50 - EnergySwap
$ws.Range("D50").Value = '9.852'
$ws.Range("E50").Value = '  -2.71%  '

# Row 51 - SynthetixNetwork
$ws.Range("D51").Value = '3.015'
$ws.Range("E51").Value = '  +0.43%  '
